# Auto-generated edit script applying the Ultima_Profits diff
# Updates currentAveragePrice / Leve price / profit columns (H-N) across
# the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR leve-profit sheets.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 1772
$ws.Range("I34").Value = 1772
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1772
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -1569
$ws.Range("N34").ClearContents()

$ws.Range("H36").Value = 1772
$ws.Range("I36").Value = 1772
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 1772
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -1057
$ws.Range("N36").ClearContents()

$ws.Range("H76").Value = 4138.7144
$ws.Range("I76").Value = 4001
$ws.Range("J76").Value = 4483
$ws.Range("K76").Value = 4001
$ws.Range("L76").Value = 4483
$ws.Range("M76").Value = -3686
$ws.Range("N76").Value = -5113

$ws.Range("H79").Value = 4138.7144
$ws.Range("I79").Value = 4001
$ws.Range("J79").Value = 4483
$ws.Range("K79").Value = 4001
$ws.Range("L79").Value = 4483
$ws.Range("M79").Value = -2909
$ws.Range("N79").Value = -6667

$ws.Range("H93").Value = 37687.75
$ws.Range("J93").Value = 37687.75
$ws.Range("L93").Value = 37687.75
$ws.Range("N93").Value = -42679.75


# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 20250
$ws.Range("I16").Value = 500
$ws.Range("J16").Value = 40000
$ws.Range("K16").Value = 500
$ws.Range("L16").Value = 40000
$ws.Range("M16").Value = -213
$ws.Range("N16").Value = -40574

$ws.Range("H32").Value = 2532.12
$ws.Range("I32").Value = 2499.111
$ws.Range("K32").Value = 2499.111
$ws.Range("M32").Value = -2212.111

$ws.Range("H45").Value = 1685151.6
$ws.Range("I45").Value = 2274338.2
$ws.Range("J45").Value = 1761.1428
$ws.Range("K45").Value = 2274338.2
$ws.Range("L45").Value = 1761.1428
$ws.Range("M45").Value = -2273961.2
$ws.Range("N45").Value = -2515.1428

$ws.Range("H61").Value = 2243.139
$ws.Range("I61").Value = 2250.0857
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 2250.0857
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -2038.0857
$ws.Range("N61").Value = -2424

$ws.Range("H117").Value = 42192.668
$ws.Range("J117").Value = 42192.668
$ws.Range("L117").Value = 42192.668
$ws.Range("N117").Value = -51370.668

$ws.Range("H122").Value = 4400.5713
$ws.Range("I122").Value = 5033.0586
$ws.Range("J122").Value = 1712.5
$ws.Range("K122").Value = 15099.1758
$ws.Range("L122").Value = 5137.5
$ws.Range("M122").Value = -12649.1758
$ws.Range("N122").Value = -10037.5

$ws.Range("H132").Value = 5436608.5
$ws.Range("I132").Value = 7813924
$ws.Range("J132").Value = 2744.7144
$ws.Range("K132").Value = 23441772
$ws.Range("L132").Value = 8234.143199999999
$ws.Range("M132").Value = -23439242
$ws.Range("N132").Value = -13294.1432

$ws.Range("H136").Value = 2243.139
$ws.Range("I136").Value = 2250.0857
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 6750.257100000001
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -4200.257100000001
$ws.Range("N136").Value = -11100


# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H51").Value = 47268
$ws.Range("J51").Value = 47268
$ws.Range("L51").Value = 47268
$ws.Range("N51").Value = -48250

$ws.Range("H52").Value = 26998.75
$ws.Range("J52").Value = 26998.75
$ws.Range("L52").Value = 26998.75
$ws.Range("N52").Value = -27524.75

$ws.Range("H121").Value = 26998.75
$ws.Range("J121").Value = 26998.75
$ws.Range("L121").Value = 26998.75
$ws.Range("N121").Value = -30492.75


# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 144.75
$ws.Range("I19").Value = 106.333336
$ws.Range("J19").Value = 260
$ws.Range("K19").Value = 106.333336
$ws.Range("L19").Value = 260
$ws.Range("M19").Value = 63.666664
$ws.Range("N19").Value = -600

$ws.Range("H23").Value = 1010
$ws.Range("J23").Value = 1010
$ws.Range("L23").Value = 1010
$ws.Range("N23").Value = -1490

$ws.Range("H24").Value = 144.75
$ws.Range("I24").Value = 106.333336
$ws.Range("J24").Value = 260
$ws.Range("K24").Value = 106.333336
$ws.Range("L24").Value = 260
$ws.Range("M24").Value = 63.666664
$ws.Range("N24").Value = -600

$ws.Range("H27").Value = 1010
$ws.Range("J27").Value = 1010
$ws.Range("L27").Value = 1010
$ws.Range("N27").Value = -1394

$ws.Range("H28").Value = 20000
$ws.Range("J28").Value = 20000
$ws.Range("L28").Value = 20000
$ws.Range("N28").Value = -20490

$ws.Range("H45").Value = 30000
$ws.Range("J45").Value = 30000
$ws.Range("L45").Value = 30000
$ws.Range("N45").Value = -31186

$ws.Range("H54").Value = 30069
$ws.Range("J54").Value = 30069
$ws.Range("L54").Value = 30069
$ws.Range("N54").Value = -31385

$ws.Range("H75").Value = 40282.57
$ws.Range("J75").Value = 40282.57
$ws.Range("L75").Value = 40282.57
$ws.Range("N75").Value = -42278.57

$ws.Range("H78").Value = 40282.57
$ws.Range("J78").Value = 40282.57
$ws.Range("L78").Value = 120847.71
$ws.Range("N78").Value = -130831.71


# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H20").Value = 3500
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()

$ws.Range("H131").Value = 1112.4648
$ws.Range("I131").Value = 632.7143
$ws.Range("J131").Value = 1230.2982
$ws.Range("K131").Value = 1898.1429
$ws.Range("L131").Value = 3690.8946
$ws.Range("M131").Value = 3141.8571
$ws.Range("N131").Value = -13770.8946


# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 4666.6665

$ws.Range("H46").Value = 19173
$ws.Range("I46").Value = 1900
$ws.Range("J46").Value = 22627.6
$ws.Range("K46").Value = 1900
$ws.Range("L46").Value = 22627.6
$ws.Range("M46").Value = -1744
$ws.Range("N46").Value = -22939.6

$ws.Range("H57").Value = 8733.166999999999
$ws.Range("J57").Value = 30000
$ws.Range("L57").Value = 30000
$ws.Range("N57").Value = -31640

$ws.Range("H70").Value = 18021.428
$ws.Range("I70").Value = 32085.715
$ws.Range("J70").Value = 3957.1428
$ws.Range("K70").Value = 32085.715
$ws.Range("L70").Value = 3957.1428
$ws.Range("M70").Value = -31815.715
$ws.Range("N70").Value = -4497.1428

$ws.Range("H73").Value = 18021.428
$ws.Range("I73").Value = 32085.715
$ws.Range("J73").Value = 3957.1428
$ws.Range("K73").Value = 32085.715
$ws.Range("L73").Value = 3957.1428
$ws.Range("M73").Value = -31149.715
$ws.Range("N73").Value = -5829.1428

$ws.Range("H80").Value = 12823194
$ws.Range("I80").Value = 22224556
$ws.Range("J80").Value = 3154.5454
$ws.Range("K80").Value = 22224556
$ws.Range("L80").Value = 3154.5454
$ws.Range("M80").Value = -22223558
$ws.Range("N80").Value = -5150.5454

$ws.Range("H83").Value = 12823194
$ws.Range("I83").Value = 22224556
$ws.Range("J83").Value = 3154.5454
$ws.Range("K83").Value = 111122780
$ws.Range("L83").Value = 15772.727
$ws.Range("M83").Value = -111117788
$ws.Range("N83").Value = -25756.727

$ws.Range("H92").Value = 23496.666
$ws.Range("J92").Value = 23496.666
$ws.Range("L92").Value = 23496.666
$ws.Range("N92").Value = -27240.666

$ws.Range("H107").Value = 994.0270400000001
$ws.Range("I107").Value = 1123
$ws.Range("J107").Value = 782.1429000000001
$ws.Range("K107").Value = 1123
$ws.Range("L107").Value = 782.1429000000001
$ws.Range("M107").Value = 797
$ws.Range("N107").Value = -4622.1429


# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 831.5833
$ws.Range("I22").Value = 318.27274
$ws.Range("J22").Value = 1265.9231
$ws.Range("K22").Value = 318.27274
$ws.Range("L22").Value = 1265.9231
$ws.Range("M22").Value = -23.27274
$ws.Range("N22").Value = -1855.9231

$ws.Range("H27").Value = 831.5833
$ws.Range("I27").Value = 318.27274
$ws.Range("J27").Value = 1265.9231
$ws.Range("K27").Value = 318.27274
$ws.Range("L27").Value = 1265.9231
$ws.Range("M27").Value = -211.27274
$ws.Range("N27").Value = -1479.9231

$ws.Range("H48").Value = 16990.166
$ws.Range("I48").Value = 18680.334
$ws.Range("J48").Value = 15300
$ws.Range("K48").Value = 18680.334
$ws.Range("L48").Value = 15300
$ws.Range("M48").Value = -18019.334
$ws.Range("N48").Value = -16622

$ws.Range("H51").Value = 23175
$ws.Range("J51").Value = 23175
$ws.Range("L51").Value = 23175
$ws.Range("N51").Value = -24131

$ws.Range("H136").Value = 3612.6538
$ws.Range("I136").Value = 1728.6809
$ws.Range("J136").Value = 21322
$ws.Range("K136").Value = 5186.0427
$ws.Range("L136").Value = 63966
$ws.Range("M136").Value = -2636.0427
$ws.Range("N136").Value = -69066


# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 2900
$ws.Range("I17").Value = 3350
$ws.Range("J17").Value = 2000
$ws.Range("K17").Value = 3350
$ws.Range("L17").Value = 2000
$ws.Range("M17").Value = -3178
$ws.Range("N17").Value = -2344

$ws.Range("H42").Value = 39424.5
$ws.Range("J42").Value = 39424.5
$ws.Range("L42").Value = 39424.5
$ws.Range("N42").Value = -40180.5

$ws.Range("H51").Value = 12018.4
$ws.Range("J51").Value = 15455
$ws.Range("L51").Value = 15455
$ws.Range("N51").Value = -16475

$ws.Range("H80").Value = 40301
$ws.Range("J80").Value = 40301
$ws.Range("L80").Value = 40301
$ws.Range("N80").Value = -42297

$ws.Range("H83").Value = 40301
$ws.Range("J83").Value = 40301
$ws.Range("L83").Value = 120903
$ws.Range("N83").Value = -130887

$ws.Range("H92").Value = 12000
$ws.Range("J92").Value = 12000
$ws.Range("L92").Value = 12000
$ws.Range("N92").Value = -16992

